$wb = $excel.ActiveWorkbook

# The CCDeferredCC_27 sheet is the closest template for the new
# CMCAutopayCC_27 test-case sheet (same headers/layout), so duplicate it to
# inherit all styles, column widths and row heights, then tweak the few
# cells that differ.
$srcSheet = $wb.Worksheets.Item("CCDeferredCC_27")
$srcSheet.Copy($null, $srcSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "CMCAutopayCC_27"

# Update the result/date/notes values for this run.
$newSheet.Range("B2").Value = "Wed Jan 29 17:30:11 IST 2025"
$newSheet.Range("E2").Value = "19"

# These two columns have no value for this test case (Amount / CalDate).
$newSheet.Range("H2").ClearContents()
$newSheet.Range("L2").ClearContents()

# Row 2 should use the default row height (no explicit custom height),
# unlike the template row that had ht="29".
$newSheet.Rows.Item(2).AutoFit()

# Match the recorded selection/active cell for the new sheet.
$newSheet.Range("L6").Select()
